# Update "想去人数" (want-to-go count) figures on both the "展览" and
# "全部类型" worksheets, which carry duplicate data.
$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    $ws.Range("F3").Value = 81
    $ws.Range("F6").Value = 30
    $ws.Range("F7").Value = 114
    $ws.Range("F8").Value = 45
}
